$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range('D2').Value = '23.783.47'
$ws.Range('E2').Value = '  -1.19%  '
$ws.Range('D3').Value = '1.636.94'
$ws.Range('E3').Value = '  -1.76%  '
$ws.Range('E4').Value = '  +0.33%  '
$ws.Range('D5').Value = '308.91'
$ws.Range('E5').Value = '  -0.54%  '
$ws.Range('E6').Value = '  +0.27%  '
$ws.Range('D7').Value = '0.3866'
$ws.Range('E7').Value = '  -1.31%  '
$ws.Range('D8').Value = '0.3803'
$ws.Range('E8').Value = '  -2.32%  '
$ws.Range('D9').Value = '50.40'
$ws.Range('E9').Value = '  -2.72%  '
$ws.Range('D10').Value = '1.318'
$ws.Range('E10').Value = '  -4.22%  '
$ws.Range('D11').Value = '1.002'
$ws.Range('E11').Value = '  +0.13%  '
$ws.Range('D12').Value = '0.08361'
$ws.Range('E12').Value = '  -1.85%  '
$ws.Range('D13').Value = '23.65'
$ws.Range('E13').Value = '  -2.36%  '
$ws.Range('D14').Value = '6.941'
$ws.Range('E14').Value = '  -4.47%  '
$ws.Range('D15').Value = '7.753'
$ws.Range('E15').Value = '  -3.60%  '
$ws.Range('D16').Value = '0.00001304'
$ws.Range('E16').Value = '  -1.42%  '
$ws.Range('D17').Value = '1.637.85'
$ws.Range('E17').Value = '  -1.19%  '
$ws.Range('D18').Value = '93.37'
$ws.Range('E18').Value = '  -1.99%  '
$ws.Range('D19').Value = '0.06930'
$ws.Range('E19').Value = '  -0.93%  '
$ws.Range('D20').Value = '19.33'
$ws.Range('E20').Value = '  -3.56%  '
$ws.Range('D21').Value = '6.844'
$ws.Range('E21').Value = '  -2.55%  '
$ws.Range('E22').Value = '  +0.30%  '
$ws.Range('D23').Value = '13.48'
$ws.Range('E23').Value = '  -2.18%  '
$ws.Range('D24').Value = '23.810.47'
$ws.Range('E24').Value = '  -1.07%  '
$ws.Range('D25').Value = '2.431'
$ws.Range('E25').Value = '  -2.28%  '
$ws.Range('D26').Value = '2.859'
$ws.Range('E26').Value = '  -10.07%  '
$ws.Range('D27').Value = '21.73'
$ws.Range('E27').Value = '  -2.72%  '
$ws.Range('D28').Value = '153.05'
$ws.Range('E28').Value = '  -0.83%  '
$ws.Range('E29').Value = '  +3.09%  '
$ws.Range('D30').Value = '136.40'
$ws.Range('E30').Value = '  -2.98%  '
$ws.Range('D31').Value = '7.796'
$ws.Range('E31').Value = '  -1.05%  '
$ws.Range('D32').Value = '2.483'
$ws.Range('E32').Value = '  -0.16%  '
$ws.Range('D33').Value = '1.818.10'
$ws.Range('E33').Value = '  -1.20%  '
$ws.Range('D34').Value = '0.07929'
$ws.Range('E34').Value = '  -3.17%  '
$ws.Range('D35').Value = '0.9781'
$ws.Range('E35').Value = '  -7.44%  '
$ws.Range('D36').Value = '0.02880'
$ws.Range('E36').Value = '  -5.19%  '
$ws.Range('D37').Value = '6.535'
$ws.Range('E37').Value = '  -3.31%  '
$ws.Range('D38').Value = '0.2643'
$ws.Range('E38').Value = '  -3.35%  '
$ws.Range('D39').Value = '10.42'
$ws.Range('E39').Value = '  -8.08%  '
$ws.Range('D40').Value = '0.09056'
$ws.Range('E40').Value = '  -1.55%  '
$ws.Range('D41').Value = '0.7440'
$ws.Range('E41').Value = '  -2.75%  '
$ws.Range('D42').Value = '1.413'
$ws.Range('E42').Value = '  -1.12%  '
$ws.Range('D43').Value = '13.19'
$ws.Range('E43').Value = '  -3.07%  '
$ws.Range('D44').Value = '16.56'
$ws.Range('E44').Value = '  -1.09%  '
$ws.Range('D45').Value = '0.6845'
$ws.Range('E45').Value = '  -3.13%  '
$ws.Range('D46').Value = '2.394'
$ws.Range('E46').Value = '  -5.00%  '
$ws.Range('D47').Value = '4.063'
$ws.Range('E47').Value = '  -0.88%  '
$ws.Range('D48').Value = '1.001'
$ws.Range('E48').Value = '  +0.27%  '
$ws.Range('D49').Value = '0.08190'
$ws.Range('E49').Value = '  -2.16%  '
$ws.Range('D50').Value = '133.65'
$ws.Range('E50').Value = '  -1.66%  '
$ws.Range('D51').Value = '1.213'
$ws.Range('E51').Value = '  -2.84%  '
